$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Cells.Item(40, 8).Value = 2248.625
$ws.Cells.Item(40, 9).Value = 1997.5
$ws.Cells.Item(40, 10).Value = 2499.75
$ws.Cells.Item(40, 11).Value = 1997.5
$ws.Cells.Item(40, 12).Value = 2499.75
$ws.Cells.Item(40, 13).Value = -1822.5
$ws.Cells.Item(40, 14).Value = -2849.75

# Row 62
$ws.Cells.Item(62, 8).Value = 7599.75
$ws.Cells.Item(62, 9).Value = 6932.6665
$ws.Cells.Item(62, 11).Value = 6932.6665
$ws.Cells.Item(62, 13).Value = -6308.6665

# Row 65
$ws.Cells.Item(65, 8).Value = 7599.75
$ws.Cells.Item(65, 9).Value = 6932.6665
$ws.Cells.Item(65, 11).Value = 34663.3325
$ws.Cells.Item(65, 13).Value = -31543.3325

# Row 106
$ws.Cells.Item(106, 8).Value = 28250.5
$ws.Cells.Item(106, 9).Value = 28250.5
$ws.Cells.Item(106, 11).Value = 28250.5
$ws.Cells.Item(106, 13).Value = -27619.5

# Row 132
$ws.Cells.Item(132, 8).Value = 1040.4412
$ws.Cells.Item(132, 9).Value = 1113.9678
$ws.Cells.Item(132, 10).Value = 280.66666
$ws.Cells.Item(132, 11).Value = 3341.9034
$ws.Cells.Item(132, 12).Value = 841.9999799999999
$ws.Cells.Item(132, 13).Value = -811.9033999999997
$ws.Cells.Item(132, 14).Value = -5901.99998

# Row 137
$ws.Cells.Item(137, 8).Value = 3420.8333
$ws.Cells.Item(137, 9).Value = 1833.1818
$ws.Cells.Item(137, 10).Value = 5915.7144
$ws.Cells.Item(137, 11).Value = 5499.5454
$ws.Cells.Item(137, 12).Value = 17747.1432
$ws.Cells.Item(137, 13).Value = -2949.5454
$ws.Cells.Item(137, 14).Value = -22847.1432

# Row 138
$ws.Cells.Item(138, 8).Value = 3454.4375
$ws.Cells.Item(138, 9).Value = 1698.8334
$ws.Cells.Item(138, 10).Value = 3636.0518
$ws.Cells.Item(138, 11).Value = 5096.5002
$ws.Cells.Item(138, 12).Value = 10908.1554
$ws.Cells.Item(138, 13).Value = 43.4997999999996
$ws.Cells.Item(138, 14).Value = -21188.1554

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 15097.728
$ws.Cells.Item(32, 10).Value = 21568.428
$ws.Cells.Item(32, 12).Value = 21568.428
$ws.Cells.Item(32, 14).Value = -22142.428

# Row 110
$ws.Cells.Item(110, 8).Value = 3213.7
$ws.Cells.Item(110, 9).Value = 3213.7
$ws.Cells.Item(110, 11).Value = 3213.7
$ws.Cells.Item(110, 13).Value = -1168.7

# Row 132
$ws.Cells.Item(132, 8).Value = 1921.409
$ws.Cells.Item(132, 10).Value = 2925.1428
$ws.Cells.Item(132, 12).Value = 8775.428400000001
$ws.Cells.Item(132, 14).Value = -13835.4284

# Row 133
$ws.Cells.Item(133, 8).Value = 49998.5
$ws.Cells.Item(133, 10).Value = 49998.5
$ws.Cells.Item(133, 12).Value = 49998.5
$ws.Cells.Item(133, 14).Value = -55058.5

# Row 140
$ws.Cells.Item(140, 8).Value = 85994.5
$ws.Cells.Item(140, 10).Value = 85994.5
$ws.Cells.Item(140, 12).Value = 85994.5
$ws.Cells.Item(140, 14).Value = -96354.5

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Cells.Item(105, 8).Value = 3199.0667
$ws.Cells.Item(105, 9).Value = 2378.6667
$ws.Cells.Item(105, 11).Value = 2378.6667
$ws.Cells.Item(105, 13).Value = -631.6667000000002

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Cells.Item(22, 8).Value = 290.25
$ws.Cells.Item(22, 9).Value = 215.1
$ws.Cells.Item(22, 10).Value = 415.5
$ws.Cells.Item(22, 11).Value = 215.1
$ws.Cells.Item(22, 12).Value = 415.5
$ws.Cells.Item(22, 13).Value = 134.9
$ws.Cells.Item(22, 14).Value = -1115.5

# Row 31
$ws.Cells.Item(31, 8).Value = 7171.3335
$ws.Cells.Item(31, 9).Value = 4000
$ws.Cells.Item(31, 10).Value = 8077.4287
$ws.Cells.Item(31, 11).Value = 4000
$ws.Cells.Item(31, 12).Value = 8077.4287
$ws.Cells.Item(31, 13).Value = -3705
$ws.Cells.Item(31, 14).Value = -8667.4287

# Row 34
$ws.Cells.Item(34, 8).Value = 7171.3335
$ws.Cells.Item(34, 9).Value = 4000
$ws.Cells.Item(34, 10).Value = 8077.4287
$ws.Cells.Item(34, 11).Value = 4000
$ws.Cells.Item(34, 12).Value = 8077.4287
$ws.Cells.Item(34, 13).Value = -3798
$ws.Cells.Item(34, 14).Value = -8481.4287

# Row 58
$ws.Cells.Item(58, 8).Value = 4226.5835
$ws.Cells.Item(58, 9).Value = 984.4
$ws.Cells.Item(58, 10).Value = 6542.4287
$ws.Cells.Item(58, 11).Value = 984.4
$ws.Cells.Item(58, 12).Value = 6542.4287
$ws.Cells.Item(58, 13).Value = -781.4
$ws.Cells.Item(58, 14).Value = -6948.4287

# Row 62
$ws.Cells.Item(62, 8).Value = 201249.5
$ws.Cells.Item(62, 10).Value = 399999
$ws.Cells.Item(62, 12).Value = 399999
$ws.Cells.Item(62, 14).Value = -401247

# Row 65
$ws.Cells.Item(65, 8).Value = 201249.5
$ws.Cells.Item(65, 10).Value = 399999
$ws.Cells.Item(65, 12).Value = 1999995
$ws.Cells.Item(65, 14).Value = -2006235

# Row 119
$ws.Cells.Item(119, 8).Value = 40000
$ws.Cells.Item(119, 10).Value = 40000
$ws.Cells.Item(119, 12).Value = 40000
$ws.Cells.Item(119, 14).Value = -49676

# Row 136
$ws.Cells.Item(136, 8).Value = 4226.5835
$ws.Cells.Item(136, 9).Value = 984.4
$ws.Cells.Item(136, 10).Value = 6542.4287
$ws.Cells.Item(136, 11).Value = 2953.2
$ws.Cells.Item(136, 12).Value = 19627.2861
$ws.Cells.Item(136, 13).Value = -403.1999999999998
$ws.Cells.Item(136, 14).Value = -24727.2861

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Cells.Item(5, 8).Value = 602.5294
$ws.Cells.Item(5, 9).Value = 643.75
$ws.Cells.Item(5, 11).Value = 1931.25
$ws.Cells.Item(5, 13).Value = -1819.25

# Row 135
$ws.Cells.Item(135, 8).Value = 602.5294
$ws.Cells.Item(135, 9).Value = 643.75
$ws.Cells.Item(135, 11).Value = 5793.75
$ws.Cells.Item(135, 13).Value = -3258.75

# Row 139
$ws.Cells.Item(139, 8).Value = 5423.8237
$ws.Cells.Item(139, 9).Value = 3213
$ws.Cells.Item(139, 11).Value = 9639
$ws.Cells.Item(139, 13).Value = -4499

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Cells.Item(80, 8).Value = 3114.84
$ws.Cells.Item(80, 9).Value = 2866.7693
$ws.Cells.Item(80, 10).Value = 3383.5833
$ws.Cells.Item(80, 11).Value = 2866.7693
$ws.Cells.Item(80, 12).Value = 3383.5833
$ws.Cells.Item(80, 13).Value = -1868.7693
$ws.Cells.Item(80, 14).Value = -5379.5833

# Row 83
$ws.Cells.Item(83, 8).Value = 3114.84
$ws.Cells.Item(83, 9).Value = 2866.7693
$ws.Cells.Item(83, 10).Value = 3383.5833
$ws.Cells.Item(83, 11).Value = 14333.8465
$ws.Cells.Item(83, 12).Value = 16917.9165
$ws.Cells.Item(83, 13).Value = -9341.8465
$ws.Cells.Item(83, 14).Value = -26901.9165

# Row 122
$ws.Cells.Item(122, 8).Value = 80268
$ws.Cells.Item(122, 9).Value = 2787.2222
$ws.Cells.Item(122, 10).Value = 254599.75
$ws.Cells.Item(122, 11).Value = 8361.6666
$ws.Cells.Item(122, 12).Value = 763799.25
$ws.Cells.Item(122, 13).Value = -5911.6666
$ws.Cells.Item(122, 14).Value = -768699.25

# Row 125
$ws.Cells.Item(125, 8).Value = 0
$ws.Cells.Item(125, 10).Value = 0
$ws.Cells.Item(125, 12).Value = 0
$ws.Cells.Item(125, 14).ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 133
$ws.Cells.Item(133, 8).Value = 119000
$ws.Cells.Item(133, 10).Value = 119000
$ws.Cells.Item(133, 12).Value = 119000
$ws.Cells.Item(133, 14).Value = -124060

$ws = $wb.Worksheets.Item("WVR")
# Row 26
$ws.Cells.Item(26, 8).Value = 1019750
$ws.Cells.Item(26, 10).Value = 2000000
$ws.Cells.Item(26, 12).Value = 2000000
$ws.Cells.Item(26, 14).Value = -2000586

# Row 29
$ws.Cells.Item(29, 8).Value = 45000
$ws.Cells.Item(29, 9).Value = 45000
$ws.Cells.Item(29, 10).Value = 0
$ws.Cells.Item(29, 11).Value = 45000
$ws.Cells.Item(29, 12).Value = 0
$ws.Cells.Item(29, 13).Value = -44710
$ws.Cells.Item(29, 14).ClearContents()

# Row 81
$ws.Cells.Item(81, 8).Value = 9599.333000000001
$ws.Cells.Item(81, 9).Value = 8320
$ws.Cells.Item(81, 10).Value = 11198.5
$ws.Cells.Item(81, 11).Value = 16640
$ws.Cells.Item(81, 12).Value = 22397
$ws.Cells.Item(81, 13).Value = -15579
$ws.Cells.Item(81, 14).Value = -24519

# Row 84
$ws.Cells.Item(84, 8).Value = 9599.333000000001
$ws.Cells.Item(84, 9).Value = 8320
$ws.Cells.Item(84, 10).Value = 11198.5
$ws.Cells.Item(84, 11).Value = 83200
$ws.Cells.Item(84, 12).Value = 111985
$ws.Cells.Item(84, 13).Value = -77896
$ws.Cells.Item(84, 14).Value = -122593

# Row 113
$ws.Cells.Item(113, 8).Value = 379.94116
$ws.Cells.Item(113, 9).Value = 492.18182
$ws.Cells.Item(113, 11).Value = 1476.54546
$ws.Cells.Item(113, 13).Value = 693.45454
